# Scheduled runner refresh: update currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ),
# and LeveProfit(NQ/HQ) cells across the Leve-profit sheets with freshly
# fetched market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 37.25
$ws.Range("I11").Value = 37.25
$ws.Range("K11").Value = 37.25
$ws.Range("M11").Value = 102.75

$ws.Range("H41").Value = 5709.65
$ws.Range("I41").Value = 501.33334
$ws.Range("J41").Value = 9971
$ws.Range("K41").Value = 501.33334
$ws.Range("L41").Value = 9971
$ws.Range("M41").Value = -61.33334000000002
$ws.Range("N41").Value = -10851

$ws.Range("H107").Value = 406.5
$ws.Range("I107").Value = 406.5
$ws.Range("K107").Value = 406.5
$ws.Range("M107").Value = 1513.5

$ws.Range("H112").Value = 1546.5143
$ws.Range("I112").Value = 1240
$ws.Range("J112").Value = 1565.091
$ws.Range("K112").Value = 3720
$ws.Range("L112").Value = 4695.272999999999
$ws.Range("M112").Value = -2612
$ws.Range("N112").Value = -6911.272999999999

$ws.Range("H132").Value = 3420.3845
$ws.Range("I132").Value = 1840
$ws.Range("J132").Value = 8161.5386
$ws.Range("K132").Value = 5520
$ws.Range("L132").Value = 24484.6158
$ws.Range("M132").Value = -2990
$ws.Range("N132").Value = -29544.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 533.3333
$ws.Range("I2").Value = 475.86206
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 475.86206
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = -362.86206
$ws.Range("N2").Value = -1176

$ws.Range("H45").Value = 3109.0908
$ws.Range("I45").Value = 1733.3334
$ws.Range("J45").Value = 3625
$ws.Range("K45").Value = 1733.3334
$ws.Range("L45").Value = 3625
$ws.Range("M45").Value = -1356.3334
$ws.Range("N45").Value = -4379

$ws.Range("H61").Value = 2435.842
$ws.Range("I61").Value = 1894.4375
$ws.Range("J61").Value = 5323.3335
$ws.Range("K61").Value = 1894.4375
$ws.Range("L61").Value = 5323.3335
$ws.Range("M61").Value = -1682.4375
$ws.Range("N61").Value = -5747.3335

$ws.Range("H102").Value = 1578.7894
$ws.Range("I102").Value = 1511.625
$ws.Range("K102").Value = 1511.625
$ws.Range("M102").Value = 110.375

$ws.Range("H116").Value = 533.3333
$ws.Range("I116").Value = 475.86206
$ws.Range("J116").Value = 950
$ws.Range("K116").Value = 475.86206
$ws.Range("L116").Value = 950
$ws.Range("M116").Value = 1818.13794
$ws.Range("N116").Value = -5538

$ws.Range("H132").Value = 4528.0884
$ws.Range("I132").Value = 5820.2
$ws.Range("J132").Value = 3508
$ws.Range("K132").Value = 17460.6
$ws.Range("L132").Value = 10524
$ws.Range("M132").Value = -14930.6
$ws.Range("N132").Value = -15584

$ws.Range("H136").Value = 2435.842
$ws.Range("I136").Value = 1894.4375
$ws.Range("J136").Value = 5323.3335
$ws.Range("K136").Value = 5683.3125
$ws.Range("L136").Value = 15970.0005
$ws.Range("M136").Value = -3133.3125
$ws.Range("N136").Value = -21070.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 533.3333
$ws.Range("I3").Value = 475.86206
$ws.Range("J3").Value = 950
$ws.Range("K3").Value = 475.86206
$ws.Range("L3").Value = 950
$ws.Range("M3").Value = -361.86206
$ws.Range("N3").Value = -1178

$ws.Range("H20").Value = 2069.7778
$ws.Range("I20").Value = 2052
$ws.Range("J20").Value = 2084
$ws.Range("K20").Value = 2052
$ws.Range("L20").Value = 2084
$ws.Range("M20").Value = -1805
$ws.Range("N20").Value = -2578

$ws.Range("H105").Value = 2431.4644
$ws.Range("I105").Value = 2144.0908
$ws.Range("J105").Value = 3485.1667
$ws.Range("K105").Value = 2144.0908
$ws.Range("L105").Value = 3485.1667
$ws.Range("M105").Value = -397.0907999999999
$ws.Range("N105").Value = -6979.1667

$ws.Range("H134").Value = 49799.13
$ws.Range("I134").Value = 63198.824
$ws.Range("J134").Value = 11833.333
$ws.Range("K134").Value = 189596.472
$ws.Range("L134").Value = 35499.999
$ws.Range("M134").Value = -187061.472
$ws.Range("N134").Value = -40569.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 352
$ws.Range("I2").Value = 352
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 352
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -239
$ws.Range("N2").ClearContents()

$ws.Range("H16").Value = 485.25
$ws.Range("I16").Value = 416.6
$ws.Range("K16").Value = 416.6
$ws.Range("M16").Value = -129.6

$ws.Range("H31").Value = 3551.2222
$ws.Range("I31").Value = 2406.1177
$ws.Range("J31").Value = 4575.7896
$ws.Range("K31").Value = 2406.1177
$ws.Range("L31").Value = 4575.7896
$ws.Range("M31").Value = -2111.1177
$ws.Range("N31").Value = -5165.7896

$ws.Range("H34").Value = 3551.2222
$ws.Range("I34").Value = 2406.1177
$ws.Range("J34").Value = 4575.7896
$ws.Range("K34").Value = 2406.1177
$ws.Range("L34").Value = 4575.7896
$ws.Range("M34").Value = -2204.1177
$ws.Range("N34").Value = -4979.7896

$ws.Range("H99").Value = 26036.143
$ws.Range("I99").Value = 45268.74
$ws.Range("J99").Value = 2754.5789
$ws.Range("K99").Value = 45268.74
$ws.Range("L99").Value = 2754.5789
$ws.Range("M99").Value = -43770.74
$ws.Range("N99").Value = -5750.5789

$ws.Range("H113").Value = 485.25
$ws.Range("I113").Value = 416.6
$ws.Range("K113").Value = 416.6
$ws.Range("M113").Value = 1753.4

$ws.Range("H126").Value = 26036.143
$ws.Range("I126").Value = 45268.74
$ws.Range("J126").Value = 2754.5789
$ws.Range("K126").Value = 135806.22
$ws.Range("L126").Value = 8263.736699999999
$ws.Range("M126").Value = -133336.22
$ws.Range("N126").Value = -13203.7367

$ws.Range("H134").Value = 2482.4443
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 11369.333
$ws.Range("I107").Value = 14460.571
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 14460.571
$ws.Range("L107").Value = 550
$ws.Range("M107").Value = -12540.571
$ws.Range("N107").Value = -4390

$ws.Range("H132").Value = 4164.1665
$ws.Range("I132").Value = 4895.778
$ws.Range("J132").Value = 3432.5557
$ws.Range("K132").Value = 14687.334
$ws.Range("L132").Value = 10297.6671
$ws.Range("M132").Value = -12157.334
$ws.Range("N132").Value = -15357.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1396.8846
$ws.Range("I16").Value = 1448.9048
$ws.Range("J16").Value = 1178.4
$ws.Range("K16").Value = 1448.9048
$ws.Range("L16").Value = 1178.4
$ws.Range("M16").Value = -1278.9048
$ws.Range("N16").Value = -1518.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2357.0715
$ws.Range("I81").Value = 1150
$ws.Range("J81").Value = 2839.9
$ws.Range("K81").Value = 2300
$ws.Range("L81").Value = 5679.8
$ws.Range("M81").Value = -1239
$ws.Range("N81").Value = -7801.8

$ws.Range("H84").Value = 2357.0715
$ws.Range("I84").Value = 1150
$ws.Range("J84").Value = 2839.9
$ws.Range("K84").Value = 11500
$ws.Range("L84").Value = 28399
$ws.Range("M84").Value = -6196
$ws.Range("N84").Value = -39007

$ws.Range("H122").Value = 47094.727
$ws.Range("I122").Value = 78414.234
$ws.Range("J122").Value = 1855.4445
$ws.Range("K122").Value = 235242.702
$ws.Range("L122").Value = 5566.333500000001
$ws.Range("M122").Value = -232792.702
$ws.Range("N122").Value = -10466.3335

$ws.Range("H132").Value = 18686.305
$ws.Range("I132").Value = 24212.232
$ws.Range("J132").Value = 3835.375
$ws.Range("K132").Value = 72636.696
$ws.Range("L132").Value = 11506.125
$ws.Range("M132").Value = -70106.696
$ws.Range("N132").Value = -16566.125

$ws.Range("H136").Value = 53336216
$ws.Range("I136").Value = 76925450
$ws.Range("J136").Value = 27781220
$ws.Range("K136").Value = 230776350
$ws.Range("L136").Value = 83343660
$ws.Range("M136").Value = -230773800
$ws.Range("N136").Value = -83348760
